$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.19849782464938
$ws.Range("C2").Value = 10.12136082065436
$ws.Range("D2").Value = 7.268754510153647
$ws.Range("F2").Value = 44.96505985538156
$ws.Range("G2").Value = 3.725170799146262
$ws.Range("I2").Value = 30.0611846307326
$ws.Range("J2").Value = 11.03391036539069
$ws.Range("K2").Value = 14.31864856071993
$ws.Range("M2").Value = 18.25461641881514
$ws.Range("N2").Value = 22.82997783874619
$ws.Range("B3").Value = 13.0205858433888
$ws.Range("C3").Value = 9.999007579124344
$ws.Range("D3").Value = 7.254562798345503
$ws.Range("F3").Value = 44.93508627768469
$ws.Range("G3").Value = 3.727946290466772
$ws.Range("I3").Value = 30.09606714905927
$ws.Range("J3").Value = 11.05015101909576
$ws.Range("K3").Value = 14.20767039524706
$ws.Range("M3").Value = 18.21320492729587
$ws.Range("N3").Value = 22.88362488930545
$ws.Range("B4").Value = 12.91374347249712
$ws.Range("C4").Value = 9.926037118394323
$ws.Range("D4").Value = 7.247123105604734
$ws.Range("F4").Value = 44.92592443785571
$ws.Range("G4").Value = 3.72974016314335
$ws.Range("I4").Value = 30.12258025204242
$ws.Range("J4").Value = 11.06159431804126
$ws.Range("K4").Value = 14.14266559405713
$ws.Range("M4").Value = 18.19144835358809
$ws.Range("N4").Value = 22.91847731930287
$ws.Range("B5").Value = 12.87086262078248
$ws.Range("C5").Value = 9.896878141385169
$ws.Range("D5").Value = 7.244414268568026
$ws.Range("F5").Value = 44.92451704690944
$ws.Range("G5").Value = 3.730493816242222
$ws.Range("I5").Value = 30.13466398265415
$ws.Range("J5").Value = 11.06662764139374
$ws.Range("K5").Value = 14.11698988299341
$ws.Range("M5").Value = 18.18351168763677
$ws.Range("N5").Value = 22.93316154830744
$ws.Range("B6").Value = 12.86378365015906
$ws.Range("C6").Value = 9.892072173381125
$ws.Range("D6").Value = 7.243984039471312
$ws.Range("F6").Value = 44.92442387805603
$ws.Range("G6").Value = 3.730620329157763
$ws.Range("I6").Value = 30.13674769361306
$ws.Range("J6").Value = 11.06748577462084
$ws.Range("K6").Value = 14.11277635025334
$ws.Range("M6").Value = 18.18225010951048
$ws.Range("N6").Value = 22.93562895636709
$ws.Range("B7").Value = 12.91316242830641
$ws.Range("C7").Value = 9.925641488103039
$ws.Range("D7").Value = 7.247085262718588
$ws.Range("F7").Value = 44.92589603719475
$ws.Range("G7").Value = 3.729750235420248
$ws.Range("I7").Value = 30.12273803973871
$ws.Range("J7").Value = 11.06166070066867
$ws.Range("K7").Value = 14.14231599275187
$ws.Range("M7").Value = 18.19133754634683
$ws.Range("N7").Value = 22.91867340533405
$ws.Range("B8").Value = 13.13669019307041
$ws.Range("C8").Value = 10.07874765838927
$ws.Range("D8").Value = 7.263598356124951
$ws.Range("F8").Value = 44.95280847034259
$ws.Range("G8").Value = 3.726109215814652
$ws.Range("I8").Value = 30.07215364384202
$ws.Range("J8").Value = 11.03920479153512
$ws.Range("K8").Value = 14.27974819761219
$ws.Range("M8").Value = 18.23958005727325
$ws.Range("N8").Value = 22.84807865129733
$ws.Range("B9").Value = 13.59154536898877
$ws.Range("C9").Value = 10.39448588215323
$ws.Range("D9").Value = 7.305969115638988
$ws.Range("F9").Value = 45.07876594242538
$ws.Range("G9").Value = 3.719677474640686
$ws.Range("I9").Value = 30.013456139805
$ws.Range("J9").Value = 11.00684057473559
$ws.Range("K9").Value = 14.57294165462707
$ws.Range("M9").Value = 18.36298930005318
$ws.Range("N9").Value = 22.72479492111473
$ws.Range("B10").Value = 13.93242387094072
$ws.Range("C10").Value = 10.63370950515409
$ws.Range("D10").Value = 7.343012614620297
$ws.Range("F10").Value = 45.21563949705784
$ws.Range("G10").Value = 3.715378940110277
$ws.Range("I10").Value = 29.99510451541554
$ws.Range("J10").Value = 10.99017215912487
$ws.Range("K10").Value = 14.80110569789359
$ws.Range("M10").Value = 18.47074592469379
$ws.Range("N10").Value = 22.64341919038837
$ws.Range("B11").Value = 14.08824018909777
$ws.Range("C11").Value = 10.74363947219159
$ws.Range("D11").Value = 7.361106610362773
$ws.Range("F11").Value = 45.28744291152618
$ws.Range("G11").Value = 3.7135150661942
$ws.Range("I11").Value = 29.9921471892626
$ws.Range("J11").Value = 10.98413107218281
$ws.Range("K11").Value = 14.90729370812736
$ws.Range("M11").Value = 18.52336065895425
$ws.Range("N11").Value = 22.60838967644954
$ws.Range("B12").Value = 14.14729314624388
$ws.Range("C12").Value = 10.78538668343437
$ws.Range("D12").Value = 7.368133108059805
$ws.Range("F12").Value = 45.3159943940189
$ws.Range("G12").Value = 3.71282235061716
$ws.Range("I12").Value = 29.99180286287973
$ws.Range("J12").Value = 10.98206487653889
$ws.Range("K12").Value = 14.94781589407213
$ws.Range("M12").Value = 18.5437904323542
$ws.Range("N12").Value = 22.59541040107795
$ws.Range("B13").Value = 14.13457379377996
$ws.Range("C13").Value = 10.77639100654553
$ws.Range("D13").Value = 7.366612122525903
$ws.Range("F13").Value = 45.30978497136218
$ws.Range("G13").Value = 3.712970958135407
$ws.Range("I13").Value = 29.9918425252397
$ws.Range("J13").Value = 10.98250002455833
$ws.Range("K13").Value = 14.9390754150298
$ws.Range("M13").Value = 18.53936821142392
$ws.Range("N13").Value = 22.59819302135771
$ws.Range("B14").Value = 14.09309786837173
$ws.Range("C14").Value = 10.74707187982371
$ws.Range("D14").Value = 7.361681207122753
$ws.Range("F14").Value = 45.28976464699563
$ws.Range("G14").Value = 3.713457814072969
$ws.Range("I14").Value = 29.99210331745611
$ws.Range("J14").Value = 10.98395664864284
$ws.Range("K14").Value = 14.91062145025697
$ws.Range("M14").Value = 18.52503136090823
$ws.Range("N14").Value = 22.60731614081258
$ws.Range("B15").Value = 14.06769728121461
$ws.Range("C15").Value = 10.72912740509051
$ws.Range("D15").Value = 7.358683506882802
$ws.Range("F15").Value = 45.27767853384231
$ws.Range("G15").Value = 3.713757730365713
$ws.Range("I15").Value = 29.99236406375421
$ws.Range("J15").Value = 10.98487770252843
$ws.Range("K15").Value = 14.89323209283292
$ws.Range("M15").Value = 18.51631513630042
$ws.Range("N15").Value = 22.61294150429509
$ws.Range("B16").Value = 13.92225090494682
$ws.Range("C16").Value = 10.62654415466912
$ws.Range("D16").Value = 7.341854798428595
$ws.Range("F16").Value = 45.21113805113312
$ws.Range("G16").Value = 3.715502584737485
$ws.Range("I16").Value = 29.99540629980312
$ws.Range("J16").Value = 10.99059795809604
$ws.Range("K16").Value = 14.79421141833537
$ws.Range("M16").Value = 18.46737882963981
$ws.Range("N16").Value = 22.64574844053983
$ws.Range("B17").Value = 13.83317369978762
$ws.Range("C17").Value = 10.5638673692236
$ws.Range("D17").Value = 7.331846437876796
$ws.Range("F17").Value = 45.17275370863779
$ws.Range("G17").Value = 3.7165963934792
$ws.Range("I17").Value = 29.99865370971684
$ws.Range("J17").Value = 10.9945018071869
$ws.Range("K17").Value = 14.73405531346263
$ws.Range("M17").Value = 18.43827109599608
$ws.Range("N17").Value = 22.66638354691676
$ws.Range("B18").Value = 13.7820128015634
$ws.Range("C18").Value = 10.52792380332108
$ws.Range("D18").Value = 7.326207118114032
$ws.Range("F18").Value = 45.15157463116884
$ws.Range("G18").Value = 3.717234144816518
$ws.Range("I18").Value = 30.00102895411402
$ws.Range("J18").Value = 10.99689230378966
$ws.Range("K18").Value = 14.69968257385431
$ws.Range("M18").Value = 18.42186826783301
$ws.Range("N18").Value = 22.67843951998079
$ws.Range("B19").Value = 13.76470512175451
$ws.Range("C19").Value = 10.51577343848948
$ws.Range("D19").Value = 7.324317993395489
$ws.Range("F19").Value = 45.14455839011159
$ws.Range("G19").Value = 3.717451559408059
$ws.Range("I19").Value = 30.00192030495166
$ws.Range("J19").Value = 10.99772661310864
$ws.Range("K19").Value = 14.6880846291522
$ws.Range("M19").Value = 18.41637313871435
$ws.Range("N19").Value = 22.68255362860839
$ws.Range("B20").Value = 13.84264889807862
$ws.Range("C20").Value = 10.57052869216117
$ws.Range("D20").Value = 7.332899741014674
$ws.Range("F20").Value = 45.17674685959232
$ws.Range("G20").Value = 3.716479063901827
$ws.Range("I20").Value = 29.99825549767528
$ws.Range("J20").Value = 10.99407121937189
$ws.Range("K20").Value = 14.74043573211275
$ws.Range("M20").Value = 18.44133464095315
$ws.Range("N20").Value = 22.6641675319751
$ws.Range("B21").Value = 14.1052794946089
$ws.Range("C21").Value = 10.75568070068802
$ws.Range("D21").Value = 7.363124828330422
$ws.Range("F21").Value = 45.29560825393211
$ws.Range("G21").Value = 3.713314457912581
$ws.Range("I21").Value = 29.99200566728944
$ws.Range("J21").Value = 10.9835227952528
$ws.Range("K21").Value = 14.91897089398764
$ws.Range("M21").Value = 18.52922881037532
$ws.Range("N21").Value = 22.60462870980736
$ws.Range("B22").Value = 14.27718132155451
$ws.Range("C22").Value = 10.87736446540734
$ws.Range("D22").Value = 7.383894917474319
$ws.Range("F22").Value = 45.38121853899768
$ws.Range("G22").Value = 3.71132248808952
$ws.Range("I22").Value = 29.99244144562447
$ws.Range("J22").Value = 10.97791933259633
$ws.Range("K22").Value = 15.03745152795059
$ws.Range("M22").Value = 18.58961467415212
$ws.Range("N22").Value = 22.56738144421222
$ws.Range("B23").Value = 14.18542962882076
$ws.Range("C23").Value = 10.81237071411121
$ws.Range("D23").Value = 7.372717920774686
$ws.Range("F23").Value = 45.33480524043911
$ws.Range("G23").Value = 3.712378683595964
$ws.Range("I23").Value = 29.99179522281539
$ws.Range("J23").Value = 10.98079200829337
$ws.Range("K23").Value = 14.97406280383927
$ws.Range("M23").Value = 18.55712034276286
$ws.Range("N23").Value = 22.58710878659102
$ws.Range("B24").Value = 13.83836499823521
$ws.Range("C24").Value = 10.5675168241873
$ws.Range("D24").Value = 7.332423185349503
$ws.Range("F24").Value = 45.17493878721429
$ws.Range("G24").Value = 3.716532080855394
$ws.Range("I24").Value = 29.99843394604883
$ws.Range("J24").Value = 10.99426543294605
$ws.Range("K24").Value = 14.73755048252459
$ws.Range("M24").Value = 18.43994857851882
$ws.Range("N24").Value = 22.66516879237145
$ws.Range("B25").Value = 13.46707500227287
$ws.Range("C25").Value = 10.30763625733711
$ws.Range("D25").Value = 7.293454547711445
$ws.Range("F25").Value = 45.03687953811198
$ws.Range("G25").Value = 3.721342114911516
$ws.Range("I25").Value = 30.02498976458906
$ws.Range("J25").Value = 11.01434685142863
$ws.Range("K25").Value = 14.49125887949696
$ws.Range("M25").Value = 18.32656561032071
$ws.Range("N25").Value = 22.7565281167587
